# Update the "trend_epi" daily RDS save schedule (row 5) so that it runs
# every day of the week (Sun-Sat => columns H:N), matching the pattern
# already used by the other rows in the schedule grid.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H5").Value = "X"
$ws.Range("J5").Value = "X"
$ws.Range("K5").Value = "X"
$ws.Range("M5").Value = "X"
$ws.Range("N5").Value = "X"

# Reflect the new active cell/selection as recorded in the workbook.
$ws.Range("N5").Select()
